{"js": "// Update the multiplication-problem table: each non-empty row of the\n// practice-sheet table has its 5 cell values replaced by new\n// \"AA\u00d7BB=\" problems, per the target revision. The new value for a given\n// cell does not depend on its old value (some values are reused\n// elsewhere in the table), so we address cells purely by their\n// (row, column) position inside the table.\n\nconst newGrid = {\n  0: [\"95\u00d730=\", \"87\u00d781=\", \"30\u00d729=\", \"89\u00d729=\", \"86\u00d721=\"],\n  4: [\"76\u00d755=\", \"72\u00d778=\", \"46\u00d721=\", \"15\u00d769=\", \"76\u00d767=\"],\n  9: [\"60\u00d752=\", \"13\u00d798=\", \"77\u00d741=\", \"90\u00d750=\", \"14\u00d717=\"],\n  14: [\"93\u00d791=\", \"15\u00d782=\", \"17\u00d715=\", \"34\u00d722=\", \"32\u00d737=\"],\n  19: [\"26\u00d774=\", \"23\u00d763=\", \"90\u00d797=\", \"71\u00d744=\", \"49\u00d776=\"],\n};\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const rowIndexStr of Object.keys(newGrid)) {\n  const rowIndex = Number(rowIndexStr);\n  const rowValues = newGrid[rowIndex];\n  for (let col = 0; col < rowValues.length; col++) {\n    table.getCell(rowIndex, col).value = rowValues[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the multiplication-problem table: each non-empty row of the\n# practice-sheet table has its 5 cell values replaced by new\n# \"AA\u00d7BB=\" problems, per the target revision. The new value for a given\n# cell does not depend on its old value (some values are reused\n# elsewhere in the table), so we address cells purely by their\n# (row, column) position inside the table (Word COM is 1-indexed).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rowMap = @{\n    1  = @(\"95\u00d730=\", \"87\u00d781=\", \"30\u00d729=\", \"89\u00d729=\", \"86\u00d721=\")\n    5  = @(\"76\u00d755=\", \"72\u00d778=\", \"46\u00d721=\", \"15\u00d769=\", \"76\u00d767=\")\n    10 = @(\"60\u00d752=\", \"13\u00d798=\", \"77\u00d741=\", \"90\u00d750=\", \"14\u00d717=\")\n    15 = @(\"93\u00d791=\", \"15\u00d782=\", \"17\u00d715=\", \"34\u00d722=\", \"32\u00d737=\")\n    20 = @(\"26\u00d774=\", \"23\u00d763=\", \"90\u00d797=\", \"71\u00d744=\", \"49\u00d776=\")\n}\n\nforeach ($rowIndex in $rowMap.Keys) {\n    $vals = $rowMap[$rowIndex]\n    for ($col = 1; $col -le $vals.Count; $col++) {\n        $t.Cell($rowIndex, $col).Range.Text = $vals[$col - 1]\n    }\n}\n"}
